# Applies the "cryptos list" update described by the commit diff.
# Cell values are written as literal text via Range.Value; numeric-looking
# price strings (column D) are apostrophe-prefixed so Excel stores them as
# text (preserving formatting such as trailing zeros) instead of coercing
# them to floating point numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '91.450.52'
$ws.Range("E2").Value = '  +3.56%  '

# Row 3
$ws.Range("D3").Value = '3.113.36'
$ws.Range("E3").Value = '  +1.23%  '

# Row 4
$ws.Range("E4").Value = '  -0.06%  '

# Row 5
$ws.Range("D5").Value = '''219.26'
$ws.Range("E5").Value = '  +4.52%  '

# Row 6
$ws.Range("D6").Value = '''621.58'
$ws.Range("E6").Value = '  -0.30%  '

# Row 7
$ws.Range("D7").Value = '''0.380'
$ws.Range("E7").Value = '  +2.16%  '

# Row 8
$ws.Range("D8").Value = '''0.893'
$ws.Range("E8").Value = '  +8.50%  '

# Row 9
$ws.Range("D9").Value = '''1.00'

# Row 10
$ws.Range("D10").Value = '3.109.69'
$ws.Range("E10").Value = '  +1.19%  '

# Row 11
$ws.Range("D11").Value = '''0.735'
$ws.Range("E11").Value = '  +23.40%  '

# Row 12
$ws.Range("D12").Value = '''0.189'
$ws.Range("E12").Value = '  +6.32%  '

# Row 13
$ws.Range("D13").Value = '''0.0000253'
$ws.Range("E13").Value = '  +6.37%  '

# Row 14
$ws.Range("D14").Value = '''5.44'
$ws.Range("E14").Value = '  +2.73%  '

# Row 15
$ws.Range("B15").Value = 'Avalanche'
$ws.Range("C15").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D15").Value = '''33.82'
$ws.Range("E15").Value = '  +7.25%  '

# Row 16
$ws.Range("B16").Value = 'WrappedBTC'
$ws.Range("C16").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D16").Value = '91.124.69'
$ws.Range("E16").Value = '  +3.41%  '

# Row 17
$ws.Range("D17").Value = '3.685.40'
$ws.Range("E17").Value = '  +1.09%  '

# Row 18
$ws.Range("D18").Value = '3.109.56'
$ws.Range("E18").Value = '  +0.89%  '

# Row 19
$ws.Range("D19").Value = '''3.77'
$ws.Range("E19").Value = '  +15.83%  '

# Row 20
$ws.Range("D20").Value = '''0.0000230'
$ws.Range("E20").Value = '  +8.92%  '

# Row 21
$ws.Range("D21").Value = '''14.06'
$ws.Range("E21").Value = '  +7.47%  '

# Row 22
$ws.Range("D22").Value = '''432.69'
$ws.Range("E22").Value = '  +2.81%  '

# Row 23
$ws.Range("D23").Value = '''8.72'
$ws.Range("E23").Value = '  +6.42%  '

# Row 24
$ws.Range("D24").Value = '''5.15'
$ws.Range("E24").Value = '  +6.68%  '

# Row 25
$ws.Range("D25").Value = '''5.93'
$ws.Range("E25").Value = '  +8.15%  '

# Row 26
$ws.Range("D26").Value = '''12.41'
$ws.Range("E26").Value = '  +8.99%  '

# Row 27
$ws.Range("D27").Value = '''83.83'
$ws.Range("E27").Value = '  +1.98%  '

# Row 28
$ws.Range("D28").Value = '3.273.87'
$ws.Range("E28").Value = '  +1.05%  '

# Row 29
$ws.Range("D29").Value = '''0.998'
$ws.Range("E29").Value = '  -0.14%  '

# Row 30
$ws.Range("D30").Value = '''0.169'
$ws.Range("E30").Value = '  +12.46%  '

# Row 31
$ws.Range("B31").Value = 'InternetComputer(DFINITY)'
$ws.Range("C31").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D31").Value = '''8.86'
$ws.Range("E31").Value = '  +10.40%  '

# Row 32
$ws.Range("B32").Value = 'Binance-PegBSC-USD'
$ws.Range("C32").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D32").Value = '''0.913'
$ws.Range("E32").Value = '  -16.10%  '

# Row 33
$ws.Range("D33").Value = '''3.86'
$ws.Range("E33").Value = '  +8.54%  '

# Row 34
$ws.Range("D34").Value = '''520.94'
$ws.Range("E34").Value = '  +3.08%  '

# Row 35
$ws.Range("D35").Value = '''7.18'
$ws.Range("E35").Value = '  +8.11%  '

# Row 36
$ws.Range("B36").Value = 'Kaspa'
$ws.Range("C36").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D36").Value = '''0.141'
$ws.Range("E36").Value = '  +8.21%  '

# Row 37
$ws.Range("B37").Value = 'Fetch.AI'
$ws.Range("C37").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D37").Value = '''1.30'
$ws.Range("E37").Value = '  +5.29%  '

# Row 38
$ws.Range("B38").Value = 'PancakeSwap'
$ws.Range("C38").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D38").Value = '''1.85'
$ws.Range("E38").Value = '  +2.50%  '

# Row 39
$ws.Range("B39").Value = 'EthereumClassic'
$ws.Range("C39").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D39").Value = '''23.24'
$ws.Range("E39").Value = '  +4.11%  '

# Row 40
$ws.Range("D40").Value = '''22.32'
$ws.Range("E40").Value = '  +0.41%  '

# Row 41
$ws.Range("E41").Value = '  -0.07%  '

# Row 42
$ws.Range("D42").Value = '''0.0775'
$ws.Range("E42").Value = '  +14.27%  '

# Row 43
$ws.Range("E43").Value = '  +0.01%  '

# Row 44
$ws.Range("B44").Value = 'Stellar'
$ws.Range("C44").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D44").Value = '''0.142'
$ws.Range("E44").Value = '  +5.82%  '

# Row 45
$ws.Range("B45").Value = 'PolygonEcosystemToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D45").Value = '''0.375'
$ws.Range("E45").Value = '  +4.58%  '

# Row 46
$ws.Range("D46").Value = '''1.90'
$ws.Range("E46").Value = '  +4.45%  '

# Row 47
$ws.Range("D47").Value = '''144.93'
$ws.Range("E47").Value = '  -1.76%  '

# Row 48
$ws.Range("D48").Value = '''44.03'
$ws.Range("E48").Value = '  +1.26%  '

# Row 49
$ws.Range("B49").Value = 'ImmutableX'
$ws.Range("C49").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D49").Value = '''1.31'
$ws.Range("E49").Value = '  +11.96%  '

# Row 50
$ws.Range("B50").Value = 'FLOKI'
$ws.Range("C50").Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range("D50").Value = '''0.000267'
$ws.Range("E50").Value = '  +24.72%  '

# Row 51
$ws.Range("D51").Value = '''166.68'
$ws.Range("E51").Value = '  +6.18%  '
